$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings are preserved literally
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.542.91"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.919.08"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "245.10"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.4851"
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.06703"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "111.29"
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("D12").Value = "1.921.77"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "0.07586"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "5.325"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.6687"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "295.84"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "30.540.95"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.534"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.173.37"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.000007528"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "6.482"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "9.432"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "164.30"
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").Value = "20.25"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "2.096"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").Value = "0.1071"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").Value = "1.444"
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("D31").Value = "4.142"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").Value = "4.040"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Value = "0.05018"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").Value = "0.7382"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "1.139"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "0.02019"
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").Value = "2.692"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "110.36"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "2.013"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "0.4408"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "0.8651"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "5.841"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "70.35"
$ws.Range("E45").Value = "  +4.81%  "
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "7.232"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "48.62"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "9.230"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "0.1225"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "0.2513"
$ws.Range("E51").Value = "  +2.87%  "
